# Generate Report for Handback
# Adds a new handback entry (7edef9fd-2185-4ad8-87dd-849db1b691a0.md) as
# row 4 on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the
# shape of the existing rows 2/3 and growing the three ListObjects
# (tables) + sheet dimensions accordingly.

$wb = $excel.ActiveWorkbook

$fileName   = "7edef9fd-2185-4ad8-87dd-849db1b691a0.md"
$pathName   = "e2e\7edef9fd-2185-4ad8-87dd-849db1b691a0.md"
$ext        = ".md"
$status     = "Handed back: in sync with en-US"
$dateFmt    = "yyyy-mm-dd HH:mm:ss"

$zhXlf      = "7edef9fd-2185-4ad8-87dd-849db1b691a0.1f55d5432e335590c0fb2f1e3ae11bc27b8936e7.zh-cn.xlf"
$zhHoDate   = "2016-08-27 08:43:57"
$zhHbDate   = "2016-08-27 08:44:27"

$deXlf      = "7edef9fd-2185-4ad8-87dd-849db1b691a0.1f55d5432e335590c0fb2f1e3ae11bc27b8936e7.de-de.xlf"
$deHoDate   = "2016-08-27 08:44:03"
$deHbDate   = "2016-08-27 08:44:34"

$srcUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7edef9fd2185a4ad887dd849db1b691a0abcdef/e2e/7edef9fd-2185-4ad8-87dd-849db1b691a0.md"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7edef9fd2185a4ad887dd849db1b691a0abcdef/e2e/7edef9fd-2185-4ad8-87dd-849db1b691a0.md"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7edef9fd2185a4ad887dd849db1b691a0abcdef/e2e/7edef9fd-2185-4ad8-87dd-849db1b691a0.md"

# ---------------------------------------------------------------------
# Sheet "Overview" (table3 / displayName "Overview") -> new row 4
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Cells.Item(4,1).Value = $fileName
$ws.Cells.Item(4,2).Value = $pathName
$ws.Cells.Item(4,3).Value = $ext
$ws.Cells.Item(4,5).Value = $status
$ws.Cells.Item(4,6).Value = $status
$ws.Cells.Item(4,7).Value = $deHoDate
$ws.Cells.Item(4,7).NumberFormat = $dateFmt

$ws.Hyperlinks.Add($ws.Cells.Item(4,2), $srcUrl, "", "", $pathName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table1) -> new row 4
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Cells.Item(4,1).Value  = $fileName
$ws.Cells.Item(4,2).Value  = $ext
$ws.Cells.Item(4,3).Value  = $status
$ws.Cells.Item(4,4).Value  = "e2e"
$ws.Cells.Item(4,5).Value  = "ht"
$ws.Cells.Item(4,6).Value  = "True"
$ws.Cells.Item(4,7).Value  = $zhXlf
$ws.Cells.Item(4,8).Value  = $zhHoDate
$ws.Cells.Item(4,8).NumberFormat = $dateFmt
$ws.Cells.Item(4,9).Value  = $fileName
$ws.Cells.Item(4,10).Value = $zhXlf
$ws.Cells.Item(4,11).Value = $zhHbDate
$ws.Cells.Item(4,11).NumberFormat = $dateFmt
$ws.Cells.Item(4,12).Value = ""
$ws.Cells.Item(4,13).Value = "True"
$ws.Cells.Item(4,14).Value = ""
$ws.Cells.Item(4,15).Value = "False"
$ws.Cells.Item(4,16).Value = ""

$ws.Hyperlinks.Add($ws.Cells.Item(4,1), $srcUrl, "", "", $fileName) | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4,9), $zhUrl, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" (table2) -> new row 4
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Cells.Item(4,1).Value  = $fileName
$ws.Cells.Item(4,2).Value  = $ext
$ws.Cells.Item(4,3).Value  = $status
$ws.Cells.Item(4,4).Value  = "e2e"
$ws.Cells.Item(4,5).Value  = "ht"
$ws.Cells.Item(4,6).Value  = "True"
$ws.Cells.Item(4,7).Value  = $deXlf
$ws.Cells.Item(4,8).Value  = $deHoDate
$ws.Cells.Item(4,8).NumberFormat = $dateFmt
$ws.Cells.Item(4,9).Value  = $fileName
$ws.Cells.Item(4,10).Value = $deXlf
$ws.Cells.Item(4,11).Value = $deHbDate
$ws.Cells.Item(4,11).NumberFormat = $dateFmt
$ws.Cells.Item(4,12).Value = ""
$ws.Cells.Item(4,13).Value = "True"
$ws.Cells.Item(4,14).Value = ""
$ws.Cells.Item(4,15).Value = "False"
$ws.Cells.Item(4,16).Value = ""

$ws.Hyperlinks.Add($ws.Cells.Item(4,1), $srcUrl, "", "", $fileName) | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4,9), $deUrl, "", "", $fileName) | Out-Null
